# Edit GroupAddr.xlsx: insert PrimaryAddrID and ActiveFlg columns before
# the existing "Query" column (which shifts from D to F), populate the
# new columns' values, and update the Query formula to include them.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank columns at D:E. This pushes the existing "Query"
# column (with its style) from D to F, preserving its original styling
# and formula (column letters inside the formula referring to B/C stay
# valid since those columns did not move).
$ws.Range("D1:E1").EntireColumn.Insert()

# The column insert shifts the sheet's trailing "default width" column
# range definition outward by two columns (e.g. max 1025 -> 1027).
# Remove the two now-superfluous trailing columns it created so the
# column range definition matches its original extent.
$ws.Columns(1026).Delete()
$ws.Columns(1026).Delete()

# New header cells
$ws.Range("D1").Value = "PrimaryAddrID"
$ws.Range("E1").Value = "ActiveFlg"

# Data for PrimaryAddrID (D) / ActiveFlg (E) per row, keyed by row number
$primaryAddrId = @{2=1; 3=0; 4=1; 5=1; 6=1; 7=0; 8=1; 9=1; 10=0; 11=1}
$activeFlg     = @{2=1; 3=1; 4=1; 5=1; 6=1; 7=1; 8=0; 9=1; 10=1; 11=1}

for ($r = 2; $r -le 11; $r++) {
    $ws.Cells.Item($r, 4).Value = $primaryAddrId[$r]
    $ws.Cells.Item($r, 5).Value = $activeFlg[$r]

    # Updated formula text referencing the relocated header row and the
    # two new value columns.
    $formula = '="INSERT INTO "&A' + $r + '&" ([" &B$1 &"],["&C$1&"],["&D$1&"],["&E$1&"]) VALUES ( ''" & B' + $r + ' & "'',''" & C' + $r + ' & "'',''" & D' + $r + ' & "'',''" & E' + $r + ' & "'')"'
    $ws.Cells.Item($r, 6).Formula = $formula
}

# Update selection to match target (F12)
$ws.Range("F12").Select()
